$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "id_test"
$ws.Range("B6").Value = "Test"
$ws.Range("A7").Value = "id_test2"
$ws.Range("B7").Value = "Test2"

$ws.Range("B10").Select() | Out-Null
